# databasestudentmodel.xlsx fixes: change_page, hasil_tes, progress, scoring, sidebars
$wb = $excel.ActiveWorkbook

$wsSiswa         = $wb.Worksheets.Item("Siswa")
$wsRiwayatKonsep = $wb.Worksheets.Item("RiwayatKonsep")
$wsRiwayatTopik  = $wb.Worksheets.Item("RiwayatTopik")
$wsKonsepTes     = $wb.Worksheets.Item("KonsepTes")
$wsTopikTes      = $wb.Worksheets.Item("TopikTes")
$wsHasilPretest  = $wb.Worksheets.Item("HasilPretest")

# --- Fix the typo in the "TopikTes" header (D1): jumah_pertanyaan -> jumlah_pertanyaan
$wsTopikTes.Range("D1").Value = "jumlah_pertanyaan"

# --- RiwayatKonsep ("RiwayatKonsep" sheet): columns resized (best-fit) to match their
# header text width, and the selection moved from D1 to C1.
$wsRiwayatKonsep.Columns.Item(2).ColumnWidth = 9.5
$wsRiwayatKonsep.Columns.Item(3).ColumnWidth = 18.166666666666668
$wsRiwayatKonsep.Columns.Item(4).ColumnWidth = 13.833333333333334
$wsRiwayatKonsep.Range("C1").Select()

# --- RiwayatTopik: columns resized, selection moved from E1 to F24
$wsRiwayatTopik.Columns.Item(2).ColumnWidth = 7.666666666666667
$wsRiwayatTopik.Columns.Item(3).ColumnWidth = 11.833333333333334
$wsRiwayatTopik.Columns.Item(4).ColumnWidth = 13.833333333333334
$wsRiwayatTopik.Columns.Item(5).ColumnWidth = 13.333333333333334
$wsRiwayatTopik.Range("F24").Select()

# --- TopikTes: columns resized, selection moved from D1 to D2, and this becomes
# the active sheet/tab (tabSelected moves here from "Siswa").
$wsTopikTes.Columns.Item(2).ColumnWidth = 5.666666666666667
$wsTopikTes.Columns.Item(4).ColumnWidth = 17
$wsTopikTes.Columns.Item(5).ColumnWidth = 13.833333333333334
$wsTopikTes.Columns.Item(6).ColumnWidth = 13.333333333333334
$wsTopikTes.Activate()
$wsTopikTes.Range("D2").Select()
